$d = $word.ActiveDocument

# --- Locate the "-Cadastramento de politica de pontuacao para descontos" paragraph ---
$findRange = $d.Content
$found = $findRange.Find.Execute("-Cadastramento de politica de pontua")
if (-not $found) {
    throw "Could not find the anchor paragraph text."
}
$anchorIndex = $findRange.Paragraphs.First.Index
$anchorPara = $d.Paragraphs.Item($anchorIndex)

# --- Split the paragraph: insert a new paragraph right after it ---
$anchorPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 1)

# Type the new line's text. A trailing sentinel character is appended and removed
# afterwards so the bookmark we (re)create below never has to be added while sitting
# at the exact "paragraph end - 1" offset (the offset right before the paragraph's
# own mark), which the COM bookmark-creation path mis-resolves.
$sentinel = "#"
$newPara.Range.Text = "-Emissão do ticket" + $sentinel

# --- Move the _GoBack bookmark from the old paragraph onto the end of the new text ---
$bookmarkPos = $newPara.Range.End - 1 - $sentinel.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- Remove the sentinel character now that the bookmark is safely placed ---
$sentinelRange = $d.Range($newPara.Range.End - 1 - $sentinel.Length, $newPara.Range.End - 1)
$sentinelRange.Text = ""

# --- Drop the now-redundant empty paragraph that used to follow the anchor paragraph ---
$trailingEmpty = $d.Paragraphs.Item($anchorIndex + 2)
$trailingEmptyLength = $trailingEmpty.Range.End - $trailingEmpty.Range.Start
if ($trailingEmptyLength -le 1) {
    $trailingEmpty.Range.Delete()
}
